$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column (H) — copy the header style from the existing "sum"
# header (G1) so H1 matches the other header cells (bold/border/centered),
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Row 2 value for the new Save column.
$ws.Range("H2").Value = 1
